$d = $word.ActiveDocument

# Locate the exact text that needs to turn into a hyperlink + extra sentence.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("ecargo@aliyun.com", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not find the target text 'ecargo@aliyun.com'"
}

# Re-materialize the found span as a fresh Range object anchored to the
# document (InsertXML on the Range object returned directly by Find does
# not reliably replace its own span in this runtime).
$rng = $d.Range($rng.Start, $rng.End)

# Run-level formatting shared by the new plain (red, bold) runs.
$rprPlain = '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'
# Run-level formatting for the hyperlink run (uses the built-in Hyperlink style, id "a4").
$rprLink = '<w:rPr><w:rStyle w:val="a4"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'

$body = ''
$body += '<w:hyperlink r:id="rIdNewMail" w:history="1"><w:r>' + $rprLink + '<w:t>ecargo@aliyun.com</w:t></w:r></w:hyperlink>'
$body += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">. </w:t></w:r>'
$body += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">Another point to note is that when uploading files, please name them with your name followed by the </w:t></w:r>'
$body += '<w:proofErr w:type="gramStart"/>'
$body += '<w:r>' + $rprPlain + '<w:t>model</w:t></w:r>'
$body += '<w:proofErr w:type="gramEnd"/>'
$body += '<w:r>' + $rprPlain + '<w:t xml:space="preserve"> name for easy distinction. For example: Qian_Jiang_GRA.docx.</w:t></w:r>'

$xmlSnippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p>' + $body + '</w:p></w:body></w:document></pkg:xmlData>' +
  '</pkg:part>' +
  '<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">' +
    '<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rIdNewMail" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="mailto:ecargo@aliyun.com" TargetMode="External"/></Relationships></pkg:xmlData>' +
  '</pkg:part>' +
'</pkg:package>'

$rng.InsertXML($xmlSnippet)
